$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 68; existing rows 68-86 shift down to 69-87.
$ws.Rows.Item(68).Insert()

# The date column (D) uses a date number format on every data row; carry that
# format over to the newly inserted row's D cell (mirrors the row below it).
$ws.Cells.Item(68, 4).NumberFormat = $ws.Cells.Item(69, 4).NumberFormat

# Populate the new row 68 with the new record's data.
$ws.Cells.Item(68, 1).Value = 7
$ws.Cells.Item(68, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(68, 3).Value = "Ñuble"
$ws.Cells.Item(68, 4).Value = 45204
$ws.Cells.Item(68, 5).Value = 16
$ws.Cells.Item(68, 6).Value = 100112026
$ws.Cells.Item(68, 7).Value = "Haba"
$ws.Cells.Item(68, 8).Value = "Sin especificar"
$ws.Cells.Item(68, 9).Value = "Primera"
$ws.Cells.Item(68, 10).Value = 50
$ws.Cells.Item(68, 11).Value = 15000
$ws.Cells.Item(68, 12).Value = 15000
$ws.Cells.Item(68, 13).Value = 15000
$ws.Cells.Item(68, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(68, 15).Value = "Región del Maule"
$ws.Cells.Item(68, 16).Value = 600
$ws.Cells.Item(68, 17).Value = 25
$ws.Cells.Item(68, 18).Value = "Hortaliza"
